# Chỉnh lại danh sách linh kiện
# Update the "VR4" potentiometer row (row 13): value 100K -> 10K, and
# give it a precise part/footprint description; the old footprint text
# ("TER-KF301-3") moves down onto the "Header 3 / Header nối biến trở"
# row (row 29), which previously had its last two columns blank; and a
# trailing spacer row (row 31) is appended below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 (VR4 biến trở): Giá trị 100K -> 10K -------------------------
$ws.Range("B13").Value = "10K"
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)   # xlPasteFormats, keep B13's own look

# --- Row 13: Kiểu chân -> part-specific description ----------------------
$ws.Range("E13").Value = "10K 5% WXD3-13-2W"
$ws.Range("E3").Copy()
$ws.Range("E13").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 29 (Header 3 / Header nối biến trở): fill in the now-vacated ---
# --- "Kiểu chân" / "Số lượng" cells with the part's old footprint data --
$ws.Range("E29").Value = "TER-KF301-3"
$ws.Range("E3").Copy()
$ws.Range("E29").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F29").Value = 1
$ws.Range("F3").Copy()
$ws.Range("F29").PasteSpecial(-4122)   # xlPasteFormats

# --- New trailing blank spacer row (row 31) under the table -------------
$a31 = $ws.Cells.Item(31, 1)
$a31.Font.Name = "Arial"
$a31.Font.Size = 9
$a31.Font.Color = 0

# --- Selection cursor, as left by the editor -----------------------------
[void]$ws.Range("I24").Select()
